$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.394592
$ws.Range("H2").Value = 1.183776
$ws.Range("I2").Value = 0.05019988976834007
$ws.Range("J2").Value = 0.05019988976834008
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 57.49234823950934
$ws.Range("R2").Value = 517.431134155584
$ws.Range("S2").Value = 0.01438697591382378
$ws.Range("T2").Value = 0.01438697591382379

$ws.Range("G3").Value = 0.394592
$ws.Range("H3").Value = 1.183776
$ws.Range("I3").Value = 0.05019988976834007
$ws.Range("J3").Value = 0.05019988976834008
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 66.60701227464534
$ws.Range("R3").Value = 599.463110471808
$ws.Range("S3").Value = 0.01666784381975465
$ws.Range("T3").Value = 0.01666784381975465

$ws.Range("G4").Value = 0.394592
$ws.Range("H4").Value = 1.183776
$ws.Range("I4").Value = 0.05019988976834007
$ws.Range("J4").Value = 0.05019988976834008
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 50.55755562222934
$ws.Range("R4").Value = 455.018000600064
$ws.Range("S4").Value = 0.0126516024700999
$ws.Range("T4").Value = 0.0126516024700999

$ws.Range("G5").Value = 0.394592
$ws.Range("H5").Value = 1.183776
$ws.Range("I5").Value = 0.05019988976834007
$ws.Range("J5").Value = 0.05019988976834008
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 25.948795684768
$ws.Range("R5").Value = 233.539161162912
$ws.Range("S5").Value = 0.006493467564661746
$ws.Range("T5").Value = 0.006493467564661747

$ws.Range("I6").Value = 0.8048623976501327
$ws.Range("J6").Value = 0.8048623976501328
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 921.783483273134
$ws.Range("R6").Value = 8296.051349458206
$ws.Range("S6").Value = 0.2306685529066215
$ws.Range("T6").Value = 0.2306685529066216

$ws.Range("I7").Value = 0.8048623976501327
$ws.Range("J7").Value = 0.8048623976501328
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.2672380517633409
$ws.Range("T7").Value = 0.2672380517633409

$ws.Range("I8").Value = 0.8048623976501327
$ws.Range("J8").Value = 0.8048623976501328
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 810.596908184864
$ws.Range("R8").Value = 7295.372173663775
$ws.Range("S8").Value = 0.2028450489670797
$ws.Range("T8").Value = 0.2028450489670797

$ws.Range("I9").Value = 0.8048623976501327
$ws.Range("J9").Value = 0.8048623976501328
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 416.0409516306869
$ws.Range("R9").Value = 3744.368564676182
$ws.Range("S9").Value = 0.1041107440130906
$ws.Range("T9").Value = 0.1041107440130907

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.078976
$ws.Range("H10").Value = 0.236928
$ws.Range("I10").Value = 0.01004730581041791
$ws.Range("J10").Value = 0.01004730581041791
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 11.50686201079467
$ws.Range("R10").Value = 103.561758097152
$ws.Range("S10").Value = 0.002879495300893447
$ws.Range("T10").Value = 0.002879495300893447

$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.078976
$ws.Range("H11").Value = 0.236928
$ws.Range("I11").Value = 0.01004730581041791
$ws.Range("J11").Value = 0.01004730581041791
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 13.33112531780267
$ws.Range("R11").Value = 119.980127860224
$ws.Range("S11").Value = 0.003336001828493591
$ws.Range("T11").Value = 0.003336001828493591

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.078976
$ws.Range("H12").Value = 0.236928
$ws.Range("I12").Value = 0.01004730581041791
$ws.Range("J12").Value = 0.01004730581041791
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 10.11889119095467
$ws.Range("R12").Value = 91.070020718592
$ws.Range("S12").Value = 0.002532167293504707
$ws.Range("T12").Value = 0.002532167293504707

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.078976
$ws.Range("H13").Value = 0.236928
$ws.Range("I13").Value = 0.01004730581041791
$ws.Range("J13").Value = 0.01004730581041791
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 5.193546975104
$ws.Range("R13").Value = 46.741922775936
$ws.Range("S13").Value = 0.001299641387526169
$ws.Range("T13").Value = 0.001299641387526169

$ws.Range("G14").Value = 1.060294666666667
$ws.Range("H14").Value = 3.180884
$ws.Range("I14").Value = 0.1348904067711093
$ws.Range("J14").Value = 0.1348904067711093
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 154.4857224994285
$ws.Range("R14").Value = 1390.371502494856
$ws.Range("S14").Value = 0.03865875088924548
$ws.Range("T14").Value = 0.03865875088924548

$ws.Range("G15").Value = 1.060294666666667
$ws.Range("H15").Value = 3.180884
$ws.Range("I15").Value = 0.1348904067711093
$ws.Range("J15").Value = 0.1348904067711093
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 178.9774244723858
$ws.Range("R15").Value = 1610.796820251472
$ws.Range("S15").Value = 0.04478759302499497
$ws.Range("T15").Value = 0.04478759302499497

$ws.Range("G16").Value = 1.060294666666667
$ws.Range("H16").Value = 3.180884
$ws.Range("I16").Value = 0.1348904067711093
$ws.Range("J16").Value = 0.1348904067711093
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 135.8514784535751
$ws.Range("R16").Value = 1222.663306082176
$ws.Range("S16").Value = 0.03399568826492617
$ws.Range("T16").Value = 0.03399568826492617

$ws.Range("G17").Value = 1.060294666666667
$ws.Range("H17").Value = 3.180884
$ws.Range("I17").Value = 0.1348904067711093
$ws.Range("J17").Value = 0.1348904067711093
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 69.72612133794534
$ws.Range("R17").Value = 627.5350920415081
$ws.Range("S17").Value = 0.01744837459194266
$ws.Range("T17").Value = 0.01744837459194266
